$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at 19 (pushes old row19 "MARIA PAULA" down to row 20) ---
$ws.Range("B19:J19").Insert(-4121)
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New row 19: re-insert JUVENAL MARTINEZ CERVANTES period 1910 (previously on row16, now row16 holds a new worker)
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1046268554"
$ws.Range("D19").Value = "JUVENAL MARTINEZ CERVANTES"
$ws.Range("E19").Value = "1910"
$ws.Range("F19").Value = 50000
$ws.Range("G19").Value = 1250000

# Row 16: replace with new worker LUIS HUMBERTO RUIZ DE LA CRUZ
$ws.Range("C16").Value = "3738908"
$ws.Range("D16").Value = "LUIS HUMBERTO RUIZ DE LA CRUZ"
$ws.Range("E16").Value = "2502"
$ws.Range("F16").Value = 70000
$ws.Range("G16").Value = 1750000

# Rows 17/18 keep JUVENAL MARTINEZ CERVANTES but the period order swaps (1912 then 1911)
$ws.Range("E17").Value = "1912"
$ws.Range("E18").Value = "1911"

# Column D needs to re-fit its width now that it holds the longer name "LUIS HUMBERTO RUIZ DE LA CRUZ"
$ws.Columns("D:D").AutoFit()

# --- Update VALOR MORA total ---
$ws.Range("E11").Value = 233200

# --- Update Cant. Trabajadores / Cant. Periodos counters ---
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 5

# NOTE: the earlier insert of row 19 already cascaded down through the rest of the sheet,
# so the old row24 ("___...") now sits at row25, and old row25 ("NOMBRE.../FIRMA...") now
# sits at row26 -- which is exactly the target layout. No further row insertion is needed.

Write-Host "done"
